$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = "28.133.11"
$ws.Range("E2").Value2 = "  -1.05%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = "1.795.18"
$ws.Range("E3").Value2 = "  +0.01%  "

$ws.Range("E4").Value2 = "  +0.00%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "1.000"
$ws.Range("E6").Value2 = "  -0.04%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = "0.5394"
$ws.Range("E7").Value2 = "  -1.17%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = "0.3783"
$ws.Range("E8").Value2 = "  -1.24%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.07451"
$ws.Range("E9").Value2 = "  -1.66%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "41.72"
$ws.Range("E10").Value2 = "  -1.65%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "1.094"
$ws.Range("E11").Value2 = "  -2.51%  "

$ws.Range("E12").Value2 = "  +0.01%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = "20.55"
$ws.Range("E13").Value2 = "  -2.58%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = "6.116"
$ws.Range("E14").Value2 = "  -1.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "7.241"
$ws.Range("E15").Value2 = "  -2.09%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = "1.789.39"
$ws.Range("E16").Value2 = "  -0.47%  "

$ws.Range("E17").Value2 = "  -2.74%  "

$ws.Range("E18").Value2 = "  -0.80%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = "0.06491"
$ws.Range("E19").Value2 = "  +0.58%  "

$ws.Range("E20").Value2 = "  -0.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "17.25"
$ws.Range("E21").Value2 = "  -0.51%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "5.901"
$ws.Range("E22").Value2 = "  -0.98%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "28.145.11"
$ws.Range("E23").Value2 = "  -0.95%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = "11.17"
$ws.Range("E24").Value2 = "  -2.07%  "

$ws.Range("E25").Value2 = "  -1.44%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "155.16"
$ws.Range("E26").Value2 = "  -2.77%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "20.29"
$ws.Range("E27").Value2 = "  -1.97%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "1.992.96"
$ws.Range("E28").Value2 = "  -0.56%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "2.282"
$ws.Range("E29").Value2 = "  -4.92%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "121.17"
$ws.Range("E30").Value2 = "  -1.53%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "1.120"
$ws.Range("E31").Value2 = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = "0.1062"
$ws.Range("E32").Value2 = "  +3.70%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.655"
$ws.Range("E33").Value2 = "  -1.08%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "5.554"
$ws.Range("E34").Value2 = "  -3.25%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "0.2256"
$ws.Range("E35").Value2 = "  -2.48%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.06489"
$ws.Range("E36").Value2 = "  +1.25%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "0.02293"
$ws.Range("E37").Value2 = "  -1.15%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = "5.014"
$ws.Range("E38").Value2 = "  -2.62%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "8.448"
$ws.Range("E39").Value2 = "  -3.69%  "

$ws.Range("B40").Value2 = "TheSandbox"
$ws.Range("C40").Value2 = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.6190"
$ws.Range("E40").Value2 = "  -3.10%  "

$ws.Range("B41").Value2 = "WEMIXTOKEN"
$ws.Range("C41").Value2 = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.450"
$ws.Range("E41").Value2 = "  +4.34%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = "11.11"
$ws.Range("E42").Value2 = "  -4.38%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "1.177"
$ws.Range("E43").Value2 = "  +1.54%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "0.9998"
$ws.Range("E44").Value2 = "  -0.10%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = "13.36"
$ws.Range("E45").Value2 = "  -1.81%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = "3.672"
$ws.Range("E46").Value2 = "  -0.03%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.5784"
$ws.Range("E47").Value2 = "  -3.01%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "124.25"
$ws.Range("E48").Value2 = "  -1.39%  "

$ws.Range("E49").Value2 = "  +3.63%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "1.926"
$ws.Range("E50").Value2 = "  -3.02%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "0.06813"
$ws.Range("E51").Value2 = "  -1.21%  "
